# Adds two new rows to the "نواقص الأصناف" (shortage items) table:
#   - "فرشه شعر اطفال الجو" inserted right after "سهايه الجو" (new row 46)
#   - "مسك الرمان" appended at the end of the items list (new row 51, just
#     before the totals row)
# Both new rows push the Total / footer rows down, the Total is
# recalculated, and the footer timestamp is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row at 46 for "فرشه شعر اطفال الجو" (copy formatting
#    from the row above it, row 45, which keeps the exact same cell
#    styles used throughout the table).
# ---------------------------------------------------------------------
$ws.Rows.Item(46).Insert()

$ws.Range("A45:Q45").Copy()
$ws.Range("A46:Q46").PasteSpecial(-4122)

$ws.Range("A46:B46").Merge()
$ws.Range("C46:G46").Merge()
$ws.Range("H46:K46").Merge()
$ws.Range("L46:M46").Merge()
$ws.Range("N46:O46").Merge()
$ws.Rows.Item(46).RowHeight = 25.5

$ws.Range("A46").Value2 = 40

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value2 = "فرشه شعر اطفال الجو"

$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value2 = "5:0"

$ws.Range("L46").NumberFormat = "@"
$ws.Range("L46").Value2 = "0"
$ws.Range("L46").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N46").NumberFormat = "@"
$ws.Range("N46").Value2 = "25.00"

$ws.Range("P46").NumberFormat = "@"
$ws.Range("P46").Value2 = "25.0000"
$ws.Range("P46").NumberFormat = "0.00"

$ws.Range("Q46").NumberFormat = "@"
$ws.Range("Q46").Value2 = "1:0"

# ---------------------------------------------------------------------
# 2) Insert a new row at 51 (right before the totals row) for
#    "مسك الرمان" - copy formatting from the row above it (row 50).
# ---------------------------------------------------------------------
$ws.Rows.Item(51).Insert()

$ws.Range("A50:Q50").Copy()
$ws.Range("A51:Q51").PasteSpecial(-4122)

$ws.Range("A51:B51").Merge()
$ws.Range("C51:G51").Merge()
$ws.Range("H51:K51").Merge()
$ws.Range("L51:M51").Merge()
$ws.Range("N51:O51").Merge()
$ws.Rows.Item(51).RowHeight = 25.5

$ws.Range("A51").Value2 = 45

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value2 = "مسك الرمان"

$ws.Range("H51").NumberFormat = "@"
$ws.Range("H51").Value2 = "3:0"

$ws.Range("L51").NumberFormat = "@"
$ws.Range("L51").Value2 = "0"
$ws.Range("L51").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N51").NumberFormat = "@"
$ws.Range("N51").Value2 = "30.00"

$ws.Range("P51").NumberFormat = "@"
$ws.Range("P51").Value2 = "30.0000"
$ws.Range("P51").NumberFormat = "0.00"

$ws.Range("Q51").NumberFormat = "@"
$ws.Range("Q51").Value2 = "1:0"

# ---------------------------------------------------------------------
# 3) Update the Total row (now row 52) with the new sum, and fix its
#    row height (25.5, matching the author's re-save).
# ---------------------------------------------------------------------
$ws.Range("P52").Value2 = 1983.28
$ws.Rows.Item(52).RowHeight = 25.5

# ---------------------------------------------------------------------
# 4) Update the footer row (now row 53) timestamp.
# ---------------------------------------------------------------------
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value2 = "Sunday, 13 July, 2025 4:02 PM"
